$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old placeholder "A." group header row (row 2). This shifts
# every row below it up by one, carrying along values/styles/text.
$ws.Rows.Item(2).Delete()

# Relabel the two "A" sub-rows (now rows 2 and 3) with their new,
# more descriptive legend text. (Assign A3 first so the new shared
# string table places "A2. ..." immediately before "A1. ...", matching
# the order the strings were authored in.)
$ws.Range("A3").Value = "A2. Lowest 3-year average:flows"
$ws.Range("A2").Value = "A1. Lowest 10-year average flows:in Reclamation's post-:2026 ensembles and:traces (2025)"

# Put the workbook's active selection back on A2 (was A10 before the edit).
$ws.Range("A2").Select()
